# Adds six new weekly-tracker rows (312-317) to Sheet1, mirroring the
# existing data layout: Participant, Date, Workout Type, Duration, Distance,
# Elevation, Zone1-5, Workout(level group), Week.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ Row=312; A="Eric";     B=45504; C="Walk";   D=18; E=0.75; F=26;  G=10; H=1;  I=0;  J=0; K=0; L="Agile Antelope"; M=8 },
    @{ Row=313; A="Steven";   B=45504; C="Walk";   D=61; E=3.57; F=673; G=73; H=8;  I=0;  J=0; K=0; L="Brave Leopard";  M=8 },
    @{ Row=314; A="Matt";     B=45504; C="Workout"; D=47; E=0;    F=0;   G=38; H=9;  I=0;  J=0; K=0; L="Agile Antelope"; M=8 },
    @{ Row=315; A="Jeremiah"; B=45505; C="Workout"; D=70; E=0;    F=0;   G=61; H=9;  I=0;  J=0; K=0; L="Agile Antelope"; M=8 },
    @{ Row=316; A="Jeremiah"; B=45505; C="Run";     D=40; E=4;    F=240; G=0;  H=18; I=15; J=1; K=0; L="Agile Antelope"; M=8 },
    @{ Row=317; A="Steven";   B=45505; C="Walk";   D=20; E=0.95; F=33;  G=20; H=0;  I=0;  J=0; K=0; L="Brave Leopard";  M=8 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    # Duplicate formatting (notably the date number format in column B) from
    # the prior row instead of assigning a NumberFormat string directly, so
    # the new cells reuse the workbook's existing style record.
    $ws.Range("A$($row - 1):M$($row - 1)").Copy($ws.Range("A$($row):M$($row)"))

    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
}

$ws.Range("A318").Select()
